$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 607, shifting existing rows 607:645 down to 608:646.
# EntireRow.Insert() copies formatting (incl. the date style on column D) from the
# row above, matching how this new weekly entry was added in the source workbook.
$ws.Rows.Item(607).Insert()

# The new row 607 is effectively a copy of the row that is now at 608 (the former
# row 607), with the weekly fields (date, volume, prices, origin) updated.
$copiedValues = $ws.Range("A608:R608").Value()
$ws.Range("A607:R607").Value = $copiedValues

# Apply the updated values for the new weekly record.
$ws.Range("D607").Value = 44783
$ws.Range("J607").Value = 480
$ws.Range("K607").Value = 36000
$ws.Range("L607").Value = 37000
$ws.Range("M607").Value = 36500
$ws.Range("O607").Value = "Provincia del Elquí"
$ws.Range("P607").Value = 521
$ws.Range("Q607").Value = 70
